$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record row (row 4) appended to the Artfynd sheet.

# Numeric columns
$ws.Range("A4").Value = 112118073
$ws.Range("B4").Value = 90658
$ws.Range("E4").Value = 4361
$ws.Range("Q4").Value = 692985.8302376649
$ws.Range("R4").Value = 6697796.684758035
$ws.Range("S4").Value = 10

# Text columns
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("F4").Value = "Orange taggsvamp"
$ws.Range("G4").Value = "Hydnellum aurantiacum"
$ws.Range("H4").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P4").Value = "Djupdal, Upl"
$ws.Range("T4").Value = "Uppsala"
$ws.Range("U4").Value = "Östhammar"
$ws.Range("V4").Value = "Uppland"
$ws.Range("W4").Value = "Gräsö"
$ws.Range("AW4").Value = "Samuel Persson"
$ws.Range("AX4").Value = "Samuel Persson"

# Date/time columns stored as literal text, not as Excel date serials.
# NumberFormat="@" forces the text-looking value to stay a string instead of
# being parsed into a date/time serial; resetting the style afterwards keeps
# the cell on the default (unformatted) style, matching the source data.
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-09-15"
$ws.Range("Y4").Style = "Normal"

$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "00:00"
$ws.Range("Z4").Style = "Normal"

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-09-15"
$ws.Range("AA4").Style = "Normal"

$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AB4").Style = "Normal"

# Boolean columns
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false

# Columns that remain blank for this record (present in the source data as
# empty string cells): I, J, K, N, AF, AT, AY
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("AF4").Value = ""
$ws.Range("AT4").Value = ""
$ws.Range("AY4").Value = ""
